# Weekly price-list update: a new daily record is inserted for
# Terminal Hortofrutícola Agro Chillán - Apio, pushing the existing
# rows 156-187 down to 157-188 (dimension grows from R187 to R188).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(156).Insert()

$ws.Cells.Item(156, 1).Value = 7
$ws.Cells.Item(156, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value = "Ñuble"
$ws.Cells.Item(156, 4).Value = 44644
$ws.Cells.Item(156, 5).Value = 16
$ws.Cells.Item(156, 6).Value = 100112017
$ws.Cells.Item(156, 7).Value = "Apio"
$ws.Cells.Item(156, 8).Value = "Americana (o)"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 60
$ws.Cells.Item(156, 11).Value = 8000
$ws.Cells.Item(156, 12).Value = 8500
$ws.Cells.Item(156, 13).Value = 8250
$ws.Cells.Item(156, 14).Value = "$/docena de matas"
$ws.Cells.Item(156, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(156, 16).Value = 1375
$ws.Cells.Item(156, 17).Value = 6
$ws.Cells.Item(156, 18).Value = "Hortaliza"
